# Weekly fruit/vegetable price update:
# Insert a new weekly record at row 97 (shifting the existing rows 97-109
# down to 98-110), matching the "Hortaliza, Agricola del Norte S.A. de
# Arica - Espinaca" weekly consolidation pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 97..109 down to 98..110, leaving row 97 free for the new entry.
$ws.Rows.Item(97).Insert()

$ws.Range("A97").Value = 1
$ws.Range("B97").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C97").Value = "Arica y Parinacota"
$ws.Range("D97").Value = 45124
$ws.Range("E97").Value = 15
$ws.Range("F97").Value = 100112012
$ws.Range("G97").Value = "Espinaca"
$ws.Range("H97").Value = "Sin especificar"
$ws.Range("I97").Value = "Segunda"
$ws.Range("J97").Value = 300
$ws.Range("K97").Value = 1300
$ws.Range("L97").Value = 1500
$ws.Range("M97").Value = 1400
$ws.Range("N97").Value = "`$/atado 2,5 a 3 kilos"
$ws.Range("O97").Value = "Región de Arica y Parinacota"
$ws.Range("P97").Value = 467
$ws.Range("Q97").Value = 3
$ws.Range("R97").Value = "Hortaliza"
